# Update "想去人数" (F column) values on both the "展览" and "全部类型"
# worksheets, which hold duplicate copies of the same data table.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1989
    4  = 122
    7  = 1658
    9  = 664
    19 = 3852
    21 = 21
    22 = 435
    24 = 712
    25 = 498
    28 = 1653
    29 = 21
    31 = 8
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
